# Fruta / hortaliza, semanal
#
# A new weekly record was inserted into the "Mango" price table at row 153,
# pushing the existing rows 153-164 down to 154-165. A couple of the
# "Origen"/price cells that follow also needed small corrections once the
# shift happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 153 - this shifts rows 153:164 down to 154:165
# and copies the row-above formatting (so D gets the date number format).
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new weekly entry.
$ws.Cells.Item(153, 1).Value = 4
$ws.Cells.Item(153, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(153, 3).Value = "Los Lagos"
$ws.Cells.Item(153, 4).Value = 44610
$ws.Cells.Item(153, 5).Value = 10
$ws.Cells.Item(153, 6).Value = "Fruta"
$ws.Cells.Item(153, 7).Value = 100108
$ws.Cells.Item(153, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(153, 9).Value = 100108002
$ws.Cells.Item(153, 10).Value = "Mango"
$ws.Cells.Item(153, 11).Value = "Sin especificar"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 120
$ws.Cells.Item(153, 14).Value = 7500
$ws.Cells.Item(153, 15).Value = 7500
$ws.Cells.Item(153, 16).Value = 7500
$ws.Cells.Item(153, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(153, 18).Value = "Perú"
$ws.Cells.Item(153, 19).Value = 1875
$ws.Cells.Item(153, 20).Value = 4

# The rows that used to be 160-163 (now 161-164) need a few price/origin
# corrections so the final figures match the source data.
$ws.Cells.Item(161, 14).Value = 7500
$ws.Cells.Item(161, 15).Value = 8000
$ws.Cells.Item(161, 16).Value = 7750
$ws.Cells.Item(161, 19).Value = 1938

$ws.Cells.Item(162, 18).Value = "Perú"
$ws.Cells.Item(162, 19).Value = 2062

$ws.Cells.Item(163, 18).Value = "Brasil"
$ws.Cells.Item(163, 19).Value = 3000
